$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.751.14"
$ws.Range("E2").Value = "  +1.17%  "

$ws.Range("D3").Value = "3.470.08"
$ws.Range("E3").Value = "  +1.25%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "414.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  -2.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.147"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.34%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.57"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.67"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.30%  "

$ws.Range("E13").Value = "  -2.39%  "

$ws.Range("D14").Value = "4.017.36"
$ws.Range("E14").Value = "  +1.23%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.140"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.26%  "

$ws.Range("D17").Value = "3.488.96"
$ws.Range("E17").Value = "  +1.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.22%  "

$ws.Range("E19").Value = "  -1.10%  "

$ws.Range("D20").Value = "62.724.66"
$ws.Range("E20").Value = "  +1.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "466.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "90.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.64%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +18.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.31%  "

$ws.Range("E28").Value = "  +0.42%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.55%  "

$ws.Range("E30").Value = "  -0.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.52%  "

$ws.Range("E32").Value = "  -2.26%  "

$ws.Range("E33").Value = "  -2.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.59%  "

$ws.Range("E35").Value = "  +0.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.51"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.54%  "

$ws.Range("E37").Value = "  -2.70%  "

$ws.Range("E38").Value = "  +4.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.06%  "

$ws.Range("E40").Value = "  -1.15%  "

$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.134"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.87%  "

$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.320"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.76%  "

$ws.Range("E43").Value = "  +7.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "146.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.90%  "

$ws.Range("E45").Value = "  +4.13%  "

$ws.Range("E46").Value = "  +1.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +11.82%  "

$ws.Range("D48").Value = "0.0₃0562"
$ws.Range("E48").Value = "  +31.65%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.83%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.140"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.60%  "
